$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the 3 new rows (116, 117, 118) ------------------------------
# Use Copy / PasteSpecial from existing rows so the new cells inherit the
# exact same style indices (s="3" date, s="4" time, s="1" wrap-text) that
# the rest of the table already uses, instead of Excel minting brand new
# numFmt/style entries.
#
# Row 112 is a good donor for rows 116-117: it has a populated "O" column
# and no "N" column, matching the target shape.
$ws.Range("A112:O112").Copy()
$ws.Range("A116:O117").PasteSpecial(-4122)
$ws.Range("N116:N117").Clear()

# Row 111 is a good donor for row 118: no "N" and no "O" column.
$ws.Range("A111:O111").Copy()
$ws.Range("A118:O118").PasteSpecial(-4122)
$ws.Range("N118:O118").Clear()

# --- Row 118 (2021-03-28, 15:13 — new unique strings created first) ----
$ws.Cells.Item(118, 5).Value = "10108170341856339"
$ws.Cells.Item(118, 4).Value = "Enjoy this weather while the mosquitoes are sluggish. Specially a mosquito. I’d like you to feel free to stay sleep."
$ws.Cells.Item(118, 1).Value = 44283
$ws.Cells.Item(118, 2).Value = 0.63472222222222219
$ws.Cells.Item(118, 3).Value = "Friends"
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = 0

# --- Row 117 (2021-03-27, 09:48) ---------------------------------------
$ws.Cells.Item(117, 5).Value = "10108167224787969"
$ws.Cells.Item(117, 4).Value = "Dense Fog bro where you at?"
$ws.Cells.Item(117, 1).Value = 44282
$ws.Cells.Item(117, 2).Value = 0.40833333333333338
$ws.Cells.Item(117, 3).Value = "Friends of Friends"
$ws.Cells.Item(117, 6).Value = 1
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 1
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 13).Value = 1
$ws.Cells.Item(117, 15).Value = "Chris Bushnell"

# --- Row 116 (2021-03-27, 09:11) ---------------------------------------
$ws.Cells.Item(116, 4).Value = "Mark McLawhorn Dense Fog Alert"
$ws.Cells.Item(116, 15).Value = "Vicky Earp"
$ws.Cells.Item(116, 5).Value = "10224539686495161"
$ws.Cells.Item(116, 1).Value = 44282
$ws.Cells.Item(116, 2).Value = 0.38263888888888892
$ws.Cells.Item(116, 3).Value = "Friends of Friends"
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 1
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 0

# --- View / selection ----------------------------------------------------
[void]$ws.Range("E122").Select()
